# Insert a new row at position 523, shifting the existing rows 523:584 down to 524:585,
# then populate the new row 523 with its data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(523).Insert()

$ws.Range("A523").Value2 = 6
$ws.Range("B523").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C523").Value2 = "Metropolitana"
$ws.Range("D523").Value2 = 44918
$ws.Range("E523").Value2 = 13
$ws.Range("F523").Value2 = 100112039
$ws.Range("G523").Value2 = "Ciboulette"
$ws.Range("H523").Value2 = "Sin especificar"
$ws.Range("I523").Value2 = "Primera"
$ws.Range("J523").Value2 = 640
$ws.Range("K523").Value2 = 700
$ws.Range("L523").Value2 = 800
$ws.Range("M523").Value2 = 744
$ws.Range("N523").Value2 = "`$/docena de atados"
$ws.Range("O523").Value2 = "Región Metropolitana"
$ws.Range("P523").Value2 = 248
$ws.Range("Q523").Value2 = 3
$ws.Range("R523").Value2 = "Hortaliza"
